# regen sval data to filter save games
# Updates the B:E (and derived G = sum) columns for rows 2-25 on the
# active worksheet with the recomputed s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B, C, D, E values per row (row number => values)
$data = @{
    2  = @(1.505614041169197,   1.65323645889881,    3.082599426703578,   0.4998867070740569)
    3  = @(3.182878228561681,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    4  = @(1.505614041169197,   0.3375848360084654,  0.7127328510149897,  0.4998867070740569)
    5  = @(3.182878228561681,   1.65323645889881,    0.7127328510149897,  0.4998867070740569)
    6  = @(1.505614041169197,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    7  = @(3.182878228561681,   1.65323645889881,    0.7127328510149897,  6.48142807727062)
    8  = @(0.3464964993005633,  0.3375848360084654,  0.7127328510149897,  0.4998867070740569)
    9  = @(0.1554434735375247,  0.3375848360084654,  0.1529057820181812,  0.4998867070740569)
    10 = @(0.06328177979961902, 1.65323645889881,    0.7127328510149897,  0.4998867070740569)
    11 = @(0.7287194209349384,  1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    12 = @(1.505614041169197,   1.65323645889881,    3.082599426703578,   0.4998867070740569)
    13 = @(3.182878228561681,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    14 = @(3.182878228561681,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    15 = @(3.182878228561681,   1.65323645889881,    16.98373111632243,   0.4998867070740569)
    16 = @(3.182878228561681,   1.65323645889881,    0.7127328510149897,  0.4998867070740569)
    17 = @(3.182878228561681,   1.65323645889881,    3.082599426703578,   0.4998867070740569)
    18 = @(0.7287194209349384,  1.65323645889881,    0.7127328510149897,  0.4998867070740569)
    19 = @(0.7287194209349384,  0.3375848360084654,  0.7127328510149897,  0.4998867070740569)
    20 = @(0.006876353814593728,0.004309184025731883,0.7127328510149897,  0.4998867070740569)
    21 = @(3.182878228561681,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    22 = @(3.182878228561681,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    23 = @(3.182878228561681,   86.29678392075563,   16.98373111632243,   6.48142807727062)
    24 = @(1.505614041169197,   1.65323645889881,    0.1529057820181812,  0.4998867070740569)
    25 = @(3.182878228561681,   1.65323645889881,    3.082599426703578,   0.4998867070740569)
}

foreach ($r in 2..25) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B = TB
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C = d2S
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D = K
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E = IP
    # G = sum, recompute from the refreshed B:E values (F/"Win" excluded)
    $ws.Cells.Item($r, 7).Value = $vals[0] + $vals[1] + $vals[2] + $vals[3]
}
